$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 40. This shifts the existing rows 40-70
# down to 41-71 (carrying their values/styles with them), matching the
# diff where every row from 41 to 71 takes on the previous row's old
# data and a new weekly entry is prepended at row 40.
$ws.Rows.Item(40).Insert()

# Populate the new row 40 with the new weekly price entry. All the
# "static" descriptive columns (B,C,E,F,G,H,I,N,O,Q,R) are identical to
# every other row in this sheet, so just repeat them; only D (date),
# K/L/M (min/max/avg price) and P (price per kg) carry genuinely new
# values per the diff (J - volume - stays 300, same as before).
$ws.Cells.Item(40, 1).Value = 1
$ws.Cells.Item(40, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(40, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(40, 4).Value = 44634
$ws.Cells.Item(40, 5).Value = 15
$ws.Cells.Item(40, 6).Value = 100112040
$ws.Cells.Item(40, 7).Value = "Cilantro"
$ws.Cells.Item(40, 8).Value = "Sin especificar"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 300
$ws.Cells.Item(40, 11).Value = 950
$ws.Cells.Item(40, 12).Value = 1000
$ws.Cells.Item(40, 13).Value = 975
$ws.Cells.Item(40, 14).Value = '$/atado 1,5 a 2 kilos'
$ws.Cells.Item(40, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(40, 16).Value = 488
$ws.Cells.Item(40, 17).Value = 2
$ws.Cells.Item(40, 18).Value = "Hortaliza"
